# hardware_map.xlsx update:
#   "more fork controls, distance sensor instead of limit"
#
# - Add a new "Distance" sensor entry in row 4 (fork mobile-goal detector)
# - Remove the "Limit*" / "lim" / "fork limit switch" entry that used to
#   live in row 26, leaving the cell blank (style retained)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4: Distance sensor used to detect mobile goals in the fork
$ws.Range("B4").Value() = "Distance"
$ws.Range("C4").Value() = "dist"
$ws.Range("D4").Value() = "detects mobile goals in fork"

# Clear out the old "Limit*" row (row 26), keeping B26's formatting/style
$ws.Range("B26:D26").ClearContents()

# Update the view: scroll so row 13 is at top, select D26
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("D26").Select()
